$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Total Cost:" label in E11, with a salmon/red fill
$ws.Range("E11").Value = "Total Cost:"
$ws.Range("E11").Interior.Color = 7962367

# Add the SUM formula in F11
$ws.Range("F11").Formula = "=SUM(F2:F10)"
$ws.Range("F11").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("F11").Font.Bold = $true

# Apply borders around the A10:H11 block (bottom-right corner of the table)
$ws.Range("A10:H10").Borders.Item(9).LineStyle = 1   # xlEdgeBottom -> thin for now; will thicken below
$ws.Range("A10:H10").Borders.Item(9).Weight = 2

$ws.Range("A11:H11").Borders.Item(7).LineStyle = 1  # left
$ws.Range("A11:H11").Borders.Item(10).LineStyle = 1 # right

# Select G19 to mirror final selection
$ws.Range("G19").Select()
